$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (class "1")
$ws.Range("B2").Value = 0.6453488372093024
$ws.Range("C2").Value = 0.7655172413793103
$ws.Range("D2").Value = 0.7003154574132492
$ws.Range("E2").Value = 145

# Row 3 (class "2")
$ws.Range("B3").Value = 0.8482758620689655
$ws.Range("C3").Value = 0.831081081081081
$ws.Range("D3").Value = 0.8395904436860068
$ws.Range("E3").Value = 148

# Row 4 (class "3")
$ws.Range("B4").Value = 0.9072847682119205
$ws.Range("C4").Value = 0.9072847682119205
$ws.Range("D4").Value = 0.9072847682119205
$ws.Range("E4").Value = 151

# Row 5 (class "4")
$ws.Range("B5").Value = 0.6742424242424242
$ws.Range("C5").Value = 0.5705128205128205
$ws.Range("D5").Value = 0.6180555555555556
$ws.Range("E5").Value = 156

# Row 6 (accuracy)
$ws.Range("B6").Value = 0.7666666666666667
$ws.Range("C6").Value = 0.7666666666666667
$ws.Range("D6").Value = 0.7666666666666667
$ws.Range("E6").Value = 0.7666666666666667

# Row 7 (macro avg)
$ws.Range("B7").Value = 0.7687879729331532
$ws.Range("C7").Value = 0.7685989777962832
$ws.Range("D7").Value = 0.7663115562166829

# Row 8 (weighted avg)
$ws.Range("B8").Value = 0.7688370452722899
$ws.Range("C8").Value = 0.7666666666666667
$ws.Range("D8").Value = 0.7653696560951947
